$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The edit is a cyclic rotation of the species-record data
# (columns A, B, E, F, G, H, Q, R) across rows 2, 3 and 4:
#   new row 2 <- old row 3
#   new row 3 <- old row 4
#   new row 4 <- old row 2
# All other columns (C, D, I..P, S..AY) stay as they are.

$cols = @("A", "B", "E", "F", "G", "H", "Q", "R")

# Capture the original values before overwriting anything.
$orig = @{}
foreach ($row in 2..4) {
    $orig[$row] = @{}
    foreach ($col in $cols) {
        $orig[$row][$col] = $ws.Range("$col$row").Value2
    }
}

$mapping = @{ 2 = 3; 3 = 4; 4 = 2 }

foreach ($destRow in 2..4) {
    $srcRow = $mapping[$destRow]
    foreach ($col in $cols) {
        $ws.Range("$col$destRow").Value2 = $orig[$srcRow][$col]
    }
}
